$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A ~56.86 chars, B = 48 chars in saved XML units) ---
$ws.Columns.Item(1).ColumnWidth = 56
$ws.Columns.Item(2).ColumnWidth = 47.166666666666664

# --- New content for the "advantages of disabling the charge controller" rows ---
# Insertion order chosen so the shared-string table fills in the same
# sequence as the source workbook.
$ws.Range("A8").Value = "4. No es necesario una LI-PO safe bag. "
$ws.Range("B7").Value = "3. disminuye consumo nominal del end-device."
$ws.Range("A7").Value = "3. Prolonga vida util de baterias LI-ION / LI-PO"
$ws.Range("A9").Value = "5. Proteccion contra corto circuito."

# --- Title row formatting ---
$title = $ws.Range("A1:B1")
$title.Font.Bold = $true
$title.Font.Size = 16
$title.HorizontalAlignment = -4108
$title.Merge()
$ws.Rows.Item(1).RowHeight = 21

# --- Center alignment for the body rows ---
$ws.Range("A4:A9").HorizontalAlignment = -4108
$ws.Range("B4:B8").HorizontalAlignment = -4108

# --- Page setup / selection ---
$ws.PageSetup.Orientation = 1
$ws.Range("A15").Select()

$wb.Save()
